$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-8.309648513793945, 25.9113883972168),
    @(-8.309648513793945, 25.9113883972168),
    @(-7.909624576568604, 27.96125221252441),
    @(-7.548133373260498, 29.9398136138916),
    @(-7.231334209442139, 31.95211791992188),
    @(-6.954930782318115, 34.16864395141602),
    @(-6.753246784210205, 36.38595581054688),
    @(-6.617630481719971, 38.65129470825195),
    @(-6.54407262802124, 40.74065399169922)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Remove the now-unused trailing rows (10-19) so the sheet's used range
# shrinks back down to A1:B9, matching the trimmed dataset.
$ws.Range("A10:B19").ClearContents()
